$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.058.36'
$ws.Range("E2").Value = '  +3.28%  '
$ws.Range("D3").Value = '2.455.59'
$ws.Range("E3").Value = '  +1.35%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '''571.88'
$ws.Range("E5").Value = '  +3.59%  '
$ws.Range("D6").Value = '''167.57'
$ws.Range("E6").Value = '  +4.30%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").Value = '''0.513'
$ws.Range("E8").Value = '  +0.55%  '
$ws.Range("E9").Value = '  +11.73%  '
$ws.Range("D10").Value = '2.453.96'
$ws.Range("E10").Value = '  +1.48%  '
$ws.Range("D11").Value = '''0.161'
$ws.Range("E11").Value = '  -1.56%  '
$ws.Range("D12").Value = '''0.336'
$ws.Range("E12").Value = '  +3.24%  '
$ws.Range("E13").Value = '  -1.37%  '
$ws.Range("D14").Value = '''0.0000182'
$ws.Range("E14").Value = '  +8.10%  '
$ws.Range("D15").Value = '69.924.49'
$ws.Range("E15").Value = '  +3.20%  '
$ws.Range("D16").Value = '2.908.66'
$ws.Range("E16").Value = '  +0.10%  '
$ws.Range("D17").Value = '''24.23'
$ws.Range("E17").Value = '  +5.32%  '
$ws.Range("D18").Value = '2.459.77'
$ws.Range("E18").Value = '  +1.41%  '
$ws.Range("D19").Value = '''10.88'
$ws.Range("E19").Value = '  +5.68%  '
$ws.Range("D20").Value = '''7.18'
$ws.Range("E20").Value = '  +5.33%  '
$ws.Range("D21").Value = '''342.91'
$ws.Range("E21").Value = '  +2.16%  '
$ws.Range("E22").Value = '  +3.26%  '
$ws.Range("D23").Value = '''2.03'
$ws.Range("E23").Value = '  +7.86%  '
$ws.Range("E24").Value = '  -0.07%  '
$ws.Range("D25").Value = '''66.58'
$ws.Range("E25").Value = '  +0.20%  '
$ws.Range("D26").Value = '''3.88'
$ws.Range("E26").Value = '  +7.29%  '
$ws.Range("B27").Value = 'WrappedeETH'
$ws.Range("C27").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D27").Value = '2.588.12'
$ws.Range("E27").Value = '  +1.51%  '
$ws.Range("B28").Value = 'Aptos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D28").Value = '''8.54'
$ws.Range("E28").Value = '  +5.68%  '
$ws.Range("D29").Value = '''0.988'
$ws.Range("E29").Value = '  -1.04%  '
$ws.Range("D30").Value = '0.0₃0860'
$ws.Range("E30").Value = '  +6.11%  '
$ws.Range("D31").Value = '''7.38'
$ws.Range("E31").Value = '  +4.01%  '
$ws.Range("E32").Value = '  +10.80%  '
$ws.Range("D33").Value = '''456.15'
$ws.Range("E33").Value = '  +8.35%  '
$ws.Range("D34").Value = '''0.999'
$ws.Range("E34").Value = '  +0.03%  '
$ws.Range("E35").Value = '  +2.60%  '
$ws.Range("D36").Value = '''160.20'
$ws.Range("E36").Value = '  -0.63%  '
$ws.Range("E37").Value = '  +9.48%  '
$ws.Range("D38").Value = '''19.10'
$ws.Range("E38").Value = '  +0.90%  '
$ws.Range("D40").Value = '''18.24'
$ws.Range("E40").Value = '  +2.68%  '
$ws.Range("D41").Value = '''0.305'
$ws.Range("D42").Value = '''1.54'
$ws.Range("E42").Value = '  +5.05%  '
$ws.Range("E43").Value = '  +4.09%  '
$ws.Range("D44").Value = '''37.98'
$ws.Range("E44").Value = '  +1.50%  '
$ws.Range("D45").Value = '''1.09'
$ws.Range("E45").Value = '  +3.17%  '
$ws.Range("D46").Value = '''2.14'
$ws.Range("E46").Value = '  +5.86%  '
$ws.Range("D47").Value = '''134.53'
$ws.Range("E47").Value = '  +4.41%  '
$ws.Range("E48").Value = '  +2.29%  '
$ws.Range("D49").Value = '''0.0729'
$ws.Range("E49").Value = '  +2.75%  '
$ws.Range("E50").Value = '  +3.02%  '
$ws.Range("D51").Value = '''0.565'
$ws.Range("E51").Value = '  +1.79%  '
